$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp_1")

# Rename the condition labels in column A of the Exp_1 sheet to more
# meaningful / readable names. The original "Cond_N" blocks are each
# 21 rows (rows 2-22, 23-43, 44-64, 65-85, 86-106).
$ws.Range("A2:A22").Value = "WT"
$ws.Range("A23:A43").Value = "KO_1"
$ws.Range("A44:A64").Value = "Another_KO"
$ws.Range("A65:A85").Value = "Overexpression"
$ws.Range("A86:A106").Value = "AnotherPertubation"

# Update the sheet selection to reflect where the edit was last made
# (the bottom condition block).
$ws.Activate() | Out-Null
$ws.Range("A86:A106").Select() | Out-Null
